# [feat] add score result to excel file
#
# The previous round's per-player score entries (B2:E4) are cleared out so
# the sheet is ready to record a new score result; the standings formulas
# in column I recompute automatically from the now-empty B:E columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old round's score numbers - the I2:I5 "25000+SUM(...)" formulas
# recalc on their own once the B:E inputs are gone.
$ws.Range("B2:E4").ClearContents()

# Re-assert the header labels (局本場 / player1..player4) so the shared
# string table is rewritten cleanly for the new round.
$ws.Range("A1").Value = "局本場"
$ws.Range("B1").Value = "player1"
$ws.Range("C1").Value = "player2"
$ws.Range("D1").Value = "player3"
$ws.Range("E1").Value = "player4"

# Re-assert the standings labels next to the running totals.
$ws.Range("H2").Value = "player1"
$ws.Range("H3").Value = "player2"
$ws.Range("H4").Value = "player3"
$ws.Range("H5").Value = "player4"

# Leave the cursor where the user would continue entering the next score.
$ws.Range("C13").Select() | Out-Null
